$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from column J (the last existing data column) to column K
$ws.Range("J3:J6").Copy() | Out-Null
$ws.Range("K3:K6").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("K3").Value = 2023
$ws.Range("K4").Value = 1462
$ws.Range("K5").Value = 462
$ws.Range("K6").Value = 1000

$ws.Range("B1:N1").ColumnWidth = 8.7109375
